$wb = $excel.ActiveWorkbook

# --- Sheet setup -----------------------------------------------------
# The workbook starts life as Sheet1 / Sheet2 / Sheet3. We want:
#   Sheet1 -> renamed to "SheetRead"  (unchanged data; just a rename + selection tweak)
#   Sheet2 -> replaced by a brand-new sheet named "SheetWrite" that is
#             populated with the "write" illustration of the same data
#   Sheet3 -> left untouched
#
# Activate the old Sheet2 first so a parameterless Worksheets.Add() inserts
# the new sheet immediately *before* it (matching Excel's default
# behaviour), which keeps the new sheet in the same tab position that
# Sheet2 used to occupy.
$wb.Worksheets.Item("Sheet2").Activate()
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "SheetWrite"

# Rename Sheet1 -> SheetRead
$wb.Worksheets.Item("Sheet1").Name = "SheetRead"

# Remove the old (now-empty, unnamed-as-Sheet2) worksheet; the new
# "SheetWrite" sheet takes its place in the tab order.
$wb.Worksheets.Item("Sheet2").Delete()

# --- SheetWrite: populate with the shopping-list data ------------------
$wsWrite = $wb.Worksheets.Item("SheetWrite")

$wsWrite.Range("A1").Value = "Id"
$wsWrite.Range("B1").Value = "Category"
$wsWrite.Range("C1").Value = "BasketItem"

$wsWrite.Range("A2").Value = 1
$wsWrite.Range("B2").Value = "Food"
$wsWrite.Range("C2").Value = "Milk"

$wsWrite.Range("A3").Value = 2
$wsWrite.Range("B3").Value = "Food"
$wsWrite.Range("C3").Value = "Bread"

$wsWrite.Range("A4").Value = 3
$wsWrite.Range("B4").Value = "Pet"
$wsWrite.Range("C4").Value = "Food"

$wsWrite.Columns.Item(2).ColumnWidth = 17.28515625
$wsWrite.Columns.Item(3).ColumnWidth = 31.7109375

# Leave SheetWrite's own cursor on C13 (per the saved file), but select it
# *before* flipping back to SheetRead below so SheetRead ends up the
# active/selected tab again (matching the original workbook, where
# Sheet1/SheetRead was -- and remains -- the tabSelected sheet).
$wsWrite.Range("C13").Select()

# --- SheetRead: just move the selection from E4 to E1, and leave it the
#     active tab -------------------------------------------------------
$wsRead = $wb.Worksheets.Item("SheetRead")
$wsRead.Activate()
$wsRead.Range("E1").Select()
